$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 984.4545
$ws.Range("I125").Value = 971.5
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 8743.5
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -6283.5
$ws.Range("N125").Value = -13920
$ws.Range("H126").Value = 23428.572
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 23428.572
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 23428.572
$ws.Range("N126").Value = -33308.572
$ws.Range("H127").Value = 83334050
$ws.Range("I127").Value = 142857660
$ws.Range("J127").Value = 980
$ws.Range("K127").Value = 428572980
$ws.Range("L127").Value = 2940
$ws.Range("M127").Value = -428568020
$ws.Range("N127").Value = -12860
$ws.Range("H128").Value = 20000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 20000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H129").Value = 894.1
$ws.Range("I129").Value = 270.66666
$ws.Range("J129").Value = 1049.9584
$ws.Range("K129").Value = 811.9999799999999
$ws.Range("L129").Value = 3149.8752
$ws.Range("M129").Value = 4188.00002
$ws.Range("N129").Value = -13149.8752
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 50001412
$ws.Range("I131").Value = 52632970
$ws.Range("J131").Value = 1800
$ws.Range("K131").Value = 157898910
$ws.Range("L131").Value = 5400
$ws.Range("M131").Value = -157893870
$ws.Range("N131").Value = -15480
$ws.Range("H132").Value = 2033.1774
$ws.Range("I132").Value = 889.2653
$ws.Range("J132").Value = 6344.846
$ws.Range("K132").Value = 2667.7959
$ws.Range("L132").Value = 19034.538
$ws.Range("M132").Value = -137.7959000000001
$ws.Range("N132").Value = -24094.538
$ws.Range("H133").Value = 56280
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 56280
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 56280
$ws.Range("N133").Value = -66400
$ws.Range("H134").Value = 31000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 31000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 31000
$ws.Range("N134").Value = -41140
$ws.Range("H135").Value = 735.08
$ws.Range("I135").Value = 724.0417
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6516.3753
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3981.3753
$ws.Range("N135").Value = -14070
$ws.Range("H136").Value = 38000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 38000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 38000
$ws.Range("N136").Value = -48200
$ws.Range("H137").Value = 1382.76
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 3397.5
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 10192.5
$ws.Range("M137").Value = -447
$ws.Range("N137").Value = -15292.5
$ws.Range("H138").Value = 2045.6
$ws.Range("I138").Value = 999.7436
$ws.Range("J138").Value = 2714.2622
$ws.Range("K138").Value = 2999.2308
$ws.Range("L138").Value = 8142.7866
$ws.Range("M138").Value = 2140.7692
$ws.Range("N138").Value = -18422.7866
$ws.Range("H139").Value = 39999.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39999.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39999.332
$ws.Range("N139").Value = -50279.332
$ws.Range("H140").Value = 40780
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 40780
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 40780
$ws.Range("N140").Value = -51140
$ws.Range("H141").Value = 2479.3333
$ws.Range("I141").Value = 2596.3845
$ws.Range("J141").Value = 2175
$ws.Range("K141").Value = 7789.1535
$ws.Range("L141").Value = 6525
$ws.Range("M141").Value = -2609.1535
$ws.Range("N141").Value = -16885

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1399.3829
$ws.Range("I122").Value = 995.6774
$ws.Range("J122").Value = 2181.5625
$ws.Range("K122").Value = 2987.0322
$ws.Range("L122").Value = 6544.6875
$ws.Range("M122").Value = -537.0322000000001
$ws.Range("N122").Value = -11444.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3399.913
$ws.Range("I132").Value = 2853.4
$ws.Range("J132").Value = 4424.625
$ws.Range("K132").Value = 8560.2
$ws.Range("L132").Value = 13273.875
$ws.Range("M132").Value = -6030.200000000001
$ws.Range("N132").Value = -18333.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4240
$ws.Range("I118").Value = 700
$ws.Range("J118").Value = 6600
$ws.Range("K118").Value = 2100
$ws.Range("L118").Value = 19800
$ws.Range("M118").Value = -857
$ws.Range("N118").Value = -22286
$ws.Range("H131").Value = 870.24
$ws.Range("I131").Value = 599.9
$ws.Range("J131").Value = 900.2778
$ws.Range("K131").Value = 1799.7
$ws.Range("L131").Value = 2700.8334
$ws.Range("M131").Value = 3240.3
$ws.Range("N131").Value = -12780.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2666.6667
$ws.Range("I80").Value = 2650
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2650
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1652
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 2666.6667
$ws.Range("I83").Value = 2650
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 13250
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -8258
$ws.Range("N83").Value = -23984
$ws.Range("H116").Value = 29932.25
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 29932.25
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 29932.25
$ws.Range("N116").Value = -39110.25
$ws.Range("H125").Value = 28494.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 28494.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 28494.5
$ws.Range("N125").Value = -33414.5
$ws.Range("H126").Value = 1534.7778
$ws.Range("I126").Value = 1458.7142
$ws.Range("J126").Value = 1801
$ws.Range("K126").Value = 4376.142599999999
$ws.Range("L126").Value = 5403
$ws.Range("M126").Value = -1906.142599999999
$ws.Range("N126").Value = -10343
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 20000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 20000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 216000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 216000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 216000
$ws.Range("N130").Value = -226040
$ws.Range("H131").Value = 28325
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 28325
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 28325
$ws.Range("N131").Value = -38405
$ws.Range("H132").Value = 2366.2144
$ws.Range("I132").Value = 2542.75
$ws.Range("J132").Value = 1801.3
$ws.Range("K132").Value = 7628.25
$ws.Range("L132").Value = 5403.9
$ws.Range("M132").Value = -5098.25
$ws.Range("N132").Value = -10463.9
$ws.Range("H133").Value = 34543
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 34543
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 34543
$ws.Range("N133").Value = -44663
$ws.Range("H134").Value = 416000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 416000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 1248000
$ws.Range("N134").Value = -1253070
$ws.Range("H135").Value = 33000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 33000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 33000
$ws.Range("N135").Value = -43140
$ws.Range("H136").Value = 19830.4
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 19830.4
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 59491.2
$ws.Range("N136").Value = -64591.2
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 47380
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47380
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47380
$ws.Range("N139").Value = -57660
$ws.Range("H140").Value = 43280
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 43280
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 43280
$ws.Range("N140").Value = -53640
$ws.Range("H141").Value = 44170.637
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 44170.637
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 44170.637
$ws.Range("N141").Value = -54530.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2191
$ws.Range("I40").Value = 2120.3
$ws.Range("J40").Value = 2426.6667
$ws.Range("K40").Value = 2120.3
$ws.Range("L40").Value = 2426.6667
$ws.Range("M40").Value = -1984.3
$ws.Range("N40").Value = -2698.6667
$ws.Range("H132").Value = 3986.44
$ws.Range("I132").Value = 4067.2
$ws.Range("J132").Value = 3865.3
$ws.Range("K132").Value = 12201.6
$ws.Range("L132").Value = 11595.9
$ws.Range("M132").Value = -9671.599999999999
$ws.Range("N132").Value = -16655.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 29900
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 29900
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29900
$ws.Range("N80").Value = -31896
$ws.Range("H83").Value = 29900
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 29900
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 89700
$ws.Range("N83").Value = -99684
$ws.Range("H113").Value = 478.125
$ws.Range("I113").Value = 372.5
$ws.Range("J113").Value = 583.75
$ws.Range("K113").Value = 1117.5
$ws.Range("L113").Value = 1751.25
$ws.Range("M113").Value = 1052.5
$ws.Range("N113").Value = -6091.25
